# Applies the "Fixed update to excel issue" commit:
#  1. Renames the "Requested quantity" header on "Weekly Quantity" to "Weekly_PO_Qty"
#  2. Renames the "Requested quantity" header on "Monthly Trend" to "Monthly_PO_Qty"
#  3. Adds a new "PO Forecast" worksheet (after "Monthly Trend") containing the
#     ds / PO_Forecast / yhat_lower / yhat_upper forecast table.

$wb = $excel.ActiveWorkbook

$wsWeekly  = $wb.Worksheets.Item("Weekly Quantity")
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")

# --- 1) Weekly Quantity header rename -------------------------------------
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

# --- 2) Monthly Trend header rename ----------------------------------------
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 3) New "PO Forecast" sheet --------------------------------------------
$sheetCount = $wb.Worksheets.Count
$wsForecast = $wb.Worksheets.Add($null, $wb.Worksheets.Item($sheetCount))
$wsForecast.Name = "PO Forecast"

# Reuse the bold/centered/bordered header style already used by the other
# sheets (style index 1 in the source workbook) ...
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)  # xlPasteFormats

# ... and the date number-format style used for column A elsewhere
# (style index 2 in the source workbook).
$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A58").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

$rows = @(
    @(2, 44969.99999999999, 171, -160.1054022116679, 525.6738136245737),
    @(3, 45025.99999999999, 184, -155.77020882793, 559.4995512241478),
    @(4, 45032.99999999999, 186, -170.8602116545419, 506.0533940785263),
    @(5, 45060.99999999999, 192, -145.9804662531609, 540.3042613891064),
    @(6, 45067.99999999999, 194, -155.0696138132409, 534.4108146741908),
    @(7, 45074.99999999999, 195, -128.8999283606942, 577.5004291084665),
    @(8, 45088.99999999999, 198, -150.8217044332972, 544.6361798811702),
    @(9, 45109.99999999999, 203, -143.0951882291968, 565.784923308461),
    @(10, 45137.99999999999, 210, -150.0596096308011, 527.5147799250241),
    @(11, 45144.99999999999, 211, -138.9772646067798, 580.454528755005),
    @(12, 45151.99999999999, 213, -157.1674481993727, 557.3187967505693),
    @(13, 45158.99999999999, 215, -141.0329074398941, 551.228757089985),
    @(14, 45172.99999999999, 218, -126.197282386975, 562.6013108578585),
    @(15, 45179.99999999999, 219, -129.768811651283, 560.6499059924888),
    @(16, 45193.99999999999, 223, -106.1199221681971, 539.9250536117619),
    @(17, 45214.99999999999, 227, -96.7400977955962, 573.0117664966532),
    @(18, 45221.99999999999, 229, -97.5716381430741, 585.3477174697647),
    @(19, 45228.99999999999, 231, -108.1749653531359, 600.7699547288375),
    @(20, 45235.99999999999, 232, -92.61705792606533, 587.5123410128617),
    @(21, 45249.99999999999, 236, -113.3161829190176, 586.776847682008),
    @(22, 45256.99999999999, 237, -113.8882346578996, 571.0615885974732),
    @(23, 45263.99999999999, 239, -95.12192127523868, 588.9100932288977),
    @(24, 45277.99999999999, 242, -121.3177597021558, 609.7162002278658),
    @(25, 45305.99999999999, 248, -119.2938802219235, 604.1281348580573),
    @(26, 45319.99999999999, 252, -93.98528371378902, 636.2666873908011),
    @(27, 45361.99999999999, 261, -91.36297529353915, 602.9307383147021),
    @(28, 45368.99999999999, 263, -101.874382054331, 593.90713531754),
    @(29, 45375.99999999999, 265, -84.18046322651443, 606.6930435930451),
    @(30, 45382.99999999999, 266, -62.39528302519009, 622.1862200133855),
    @(31, 45389.99999999999, 268, -115.0060939533052, 613.9378168416059),
    @(32, 45396.99999999999, 269, -99.97918140564859, 598.1289717388545),
    @(33, 45403.99999999999, 271, -50.75941951377362, 637.8824467883907),
    @(34, 45410.99999999999, 273, -90.07127640427646, 607.1563885752696),
    @(35, 45424.99999999999, 276, -64.8234871109966, 594.6125011023346),
    @(36, 45431.99999999999, 277, -62.1576247861527, 633.8885518651692),
    @(37, 45438.99999999999, 279, -74.87781931351586, 627.7641181791848),
    @(38, 45445.99999999999, 281, -75.63591172742719, 631.7937822721158),
    @(39, 45459.99999999999, 284, -55.64907577599489, 634.6292188570059),
    @(40, 45466.99999999999, 285, -64.62295862970149, 637.7926876115358),
    @(41, 45487.99999999999, 290, -48.80002005362939, 627.1318900110194),
    @(42, 45494.99999999999, 292, -44.74113161221559, 625.1002197896214),
    @(43, 45501.99999999999, 294, -89.32003917600154, 639.0209036286878),
    @(44, 45508.99999999999, 295, -70.30396197909454, 654.3385663237872),
    @(45, 45515.99999999999, 297, -28.25909431435525, 659.1427863372119),
    @(46, 45522.99999999999, 298, -42.07687703276896, 631.5800850529752),
    @(47, 45529.99999999999, 300, -52.58861928669527, 636.8393044544423),
    @(48, 45557.99999999999, 306, -44.33087997744075, 630.9504537356979),
    @(49, 45564.99999999999, 308, -52.45066551769133, 646.7829167686053),
    @(50, 45571.99999999999, 310, -31.78897189376545, 657.2490376909783),
    @(51, 45578.99999999999, 311, -17.92426633937998, 676.4969254458767),
    @(52, 45585.99999999999, 313, -32.67790171500803, 640.368797115256),
    @(53, 45592.99999999999, 314, -32.21221043642084, 686.9172605297214),
    @(54, 45599.99999999999, 316, -39.10630610469349, 637.6659599390197),
    @(55, 45606.99999999999, 318, -30.89494716277663, 651.3209284834546),
    @(56, 45613.99999999999, 319, -14.08062008535468, 666.8265480945553),
    @(57, 45620.99999999999, 321, -17.78763219181065, 673.0828575459514),
    @(58, 45627.99999999999, 323, -14.92607638750785, 660.4660959904891)
)

foreach ($row in $rows) {
    $r = $row[0]
    $wsForecast.Cells.Item($r, 1).Value = $row[1]
    $wsForecast.Cells.Item($r, 2).Value = $row[2]
    $wsForecast.Cells.Item($r, 3).Value = $row[3]
    $wsForecast.Cells.Item($r, 4).Value = $row[4]
}

# Restore the original active sheet/selection (adding a sheet makes it active).
$wsWeekly.Activate() | Out-Null
$wsWeekly.Range("A1").Select() | Out-Null
